# Applies the cryptos list price/volume update described by the commit
# "Updated cryptos list on Fri Nov  3 10:30:15 UTC 2023 with GitHub Actions".
# Values are written with a leading apostrophe so Excel stores them as literal
# text (matching the original inlineStr cells) instead of re-parsing strings
# like "34.553.38" or "11.12" as numbers/dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "'34.553.38" },
    @{ Cell = "E2"; Value = "'  -2.88%  " },
    @{ Cell = "D3"; Value = "'1.801.76" },
    @{ Cell = "E3"; Value = "'  -2.23%  " },
    @{ Cell = "E4"; Value = "'  +0.46%  " },
    @{ Cell = "D5"; Value = "'228.76" },
    @{ Cell = "E5"; Value = "'  -1.38%  " },
    @{ Cell = "E6"; Value = "'  -1.31%  " },
    @{ Cell = "E7"; Value = "'  +0.56%  " },
    @{ Cell = "D8"; Value = "'38.90" },
    @{ Cell = "E8"; Value = "'  -11.29%  " },
    @{ Cell = "D9"; Value = "'0.320" },
    @{ Cell = "E9"; Value = "'  +2.69%  " },
    @{ Cell = "D10"; Value = "'0.0677" },
    @{ Cell = "E10"; Value = "'  -4.09%  " },
    @{ Cell = "E11"; Value = "'  -2.12%  " },
    @{ Cell = "D12"; Value = "'2.062.94" },
    @{ Cell = "E12"; Value = "'  -2.20%  " },
    @{ Cell = "B13"; Value = "'Chainlink" },
    @{ Cell = "C13"; Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link" },
    @{ Cell = "D13"; Value = "'11.12" },
    @{ Cell = "E13"; Value = "'  -1.66%  " },
    @{ Cell = "B14"; Value = "'WrappedEther" },
    @{ Cell = "C14"; Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth" },
    @{ Cell = "D14"; Value = "'1.810.90" },
    @{ Cell = "E14"; Value = "'  -1.83%  " },
    @{ Cell = "D15"; Value = "'0.658" },
    @{ Cell = "E15"; Value = "'  -2.67%  " },
    @{ Cell = "D16"; Value = "'4.55" },
    @{ Cell = "E16"; Value = "'  -4.36%  " },
    @{ Cell = "D17"; Value = "'34.611.75" },
    @{ Cell = "E17"; Value = "'  -2.60%  " },
    @{ Cell = "D18"; Value = "'68.90" },
    @{ Cell = "E18"; Value = "'  -2.26%  " },
    @{ Cell = "D19"; Value = "'242.83" },
    @{ Cell = "E19"; Value = "'  -1.07%  " },
    @{ Cell = "D20"; Value = "'0.0₃0777" },
    @{ Cell = "E20"; Value = "'  -3.38%  " },
    @{ Cell = "D21"; Value = "'11.74" },
    @{ Cell = "E21"; Value = "'  -2.67%  " },
    @{ Cell = "E22"; Value = "'  -1.32%  " },
    @{ Cell = "E23"; Value = "'  +0.54%  " },
    @{ Cell = "D24"; Value = "'2.23" },
    @{ Cell = "E24"; Value = "'  +0.19%  " },
    @{ Cell = "D25"; Value = "'171.91" },
    @{ Cell = "E25"; Value = "'  -0.27%  " },
    @{ Cell = "D26"; Value = "'7.70" },
    @{ Cell = "D27"; Value = "'17.11" },
    @{ Cell = "E27"; Value = "'  -4.28%  " },
    @{ Cell = "E28"; Value = "'  -1.44%  " },
    @{ Cell = "D29"; Value = "'1.48" },
    @{ Cell = "E29"; Value = "'  -5.55%  " },
    @{ Cell = "E30"; Value = "'  +0.51%  " },
    @{ Cell = "D31"; Value = "'4.03" },
    @{ Cell = "E31"; Value = "'  +1.90%  " },
    @{ Cell = "D32"; Value = "'0.0539" },
    @{ Cell = "E32"; Value = "'  -2.50%  " },
    @{ Cell = "D33"; Value = "'3.87" },
    @{ Cell = "E33"; Value = "'  -5.32%  " },
    @{ Cell = "D34"; Value = "'1.22" },
    @{ Cell = "E34"; Value = "'  +6.42%  " },
    @{ Cell = "E35"; Value = "'  -4.14%  " },
    @{ Cell = "E36"; Value = "'  -0.54%  " },
    @{ Cell = "D37"; Value = "'90.83" },
    @{ Cell = "E37"; Value = "'  -5.50%  " },
    @{ Cell = "E38"; Value = "'  +4.11%  " },
    @{ Cell = "D39"; Value = "'1.313.20" },
    @{ Cell = "E39"; Value = "'  -2.86%  " },
    @{ Cell = "E40"; Value = "'  -2.53%  " },
    @{ Cell = "E41"; Value = "'  -0.37%  " },
    @{ Cell = "D42"; Value = "'0.952" },
    @{ Cell = "E42"; Value = "'  -6.54%  " },
    @{ Cell = "D43"; Value = "'14.28" },
    @{ Cell = "E43"; Value = "'  -7.82%  " },
    @{ Cell = "B44"; Value = "'RenderToken" },
    @{ Cell = "C44"; Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr" },
    @{ Cell = "D44"; Value = "'2.19" },
    @{ Cell = "E44"; Value = "'  -11.20%  " },
    @{ Cell = "B45"; Value = "'MXToken" },
    @{ Cell = "C45"; Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx" },
    @{ Cell = "D45"; Value = "'2.70" },
    @{ Cell = "E45"; Value = "'  -4.07%  " },
    @{ Cell = "D46"; Value = "'6.17" },
    @{ Cell = "E46"; Value = "'  -2.14%  " },
    @{ Cell = "D47"; Value = "'0.0512" },
    @{ Cell = "E47"; Value = "'  -1.12%  " },
    @{ Cell = "D48"; Value = "'1.984.57" },
    @{ Cell = "E48"; Value = "'  -1.40%  " },
    @{ Cell = "E49"; Value = "'  +0.55%  " },
    @{ Cell = "E50"; Value = "'  +4.05%  " },
    @{ Cell = "D51"; Value = "'97.48" },
    @{ Cell = "E51"; Value = "'  -5.38%  " }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
